$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff effectively swaps the weekly price-record values between row 2
# and row 4 (columns D, M, N, O, P, R, S), while leaving everything else
# (A, B, C, E-L, Q, T) untouched. Apply the exact target values directly.

# Row 2 (was: 2022-06-07 / Limarí record) becomes the 2021-06-15 / Curicó record
$ws.Range("D2").Value = 44362
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1083

# Row 4 (was: 2021-06-15 / Curicó record) becomes the 2022-06-07 / Limarí record
$ws.Range("D4").Value = 44719
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20400
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1133
